# Updates cryptos list values per the commit "Updated cryptos list on Wed Jul 12 15:36:38 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text storage so numeric-looking strings (e.g. "246.16") are not
    # silently reinterpreted as numbers by Excel, then restore the default
    # style so no stray number-format style is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 (D,E)
$ws.Cells.Item(2, 4).Value = '30.681.71'
$ws.Cells.Item(2, 5).Value = '  +0.35%  '

# Row 3 (D,E)
$ws.Cells.Item(3, 4).Value = '1.896.13'
$ws.Cells.Item(3, 5).Value = '  +1.28%  '

# Row 4 (E)
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5 (D,E)
Set-TextCell 5 4 '246.16'
$ws.Cells.Item(5, 5).Value = '  -0.25%  '

# Row 6 (E)
$ws.Cells.Item(6, 5).Value = '  -0.01%  '

# Row 7 (D,E)
Set-TextCell 7 4 '0.4727'
$ws.Cells.Item(7, 5).Value = '  -0.22%  '

# Row 8 (D,E)
Set-TextCell 8 4 '0.2925'
$ws.Cells.Item(8, 5).Value = '  +0.59%  '

# Row 9 (D,E)
Set-TextCell 9 4 '0.06516'
$ws.Cells.Item(9, 5).Value = '  +0.57%  '

# Row 10 (D,E)
Set-TextCell 10 4 '22.58'
$ws.Cells.Item(10, 5).Value = '  +2.81%  '

# Row 11 (D,E)
Set-TextCell 11 4 '0.07780'
$ws.Cells.Item(11, 5).Value = '  +0.59%  '

# Row 12 (D,E)
$ws.Cells.Item(12, 4).Value = '1.893.28'
$ws.Cells.Item(12, 5).Value = '  +1.28%  '

# Row 13 (E)
$ws.Cells.Item(13, 5).Value = '  +0.04%  '

# Row 14 (D,E)
Set-TextCell 14 4 '96.61'
$ws.Cells.Item(14, 5).Value = '  +0.13%  '

# Row 15 (D,E)
Set-TextCell 15 4 '5.216'
$ws.Cells.Item(15, 5).Value = '  +1.26%  '

# Row 16 (D,E)
Set-TextCell 16 4 '284.88'
$ws.Cells.Item(16, 5).Value = '  +4.01%  '

# Row 17 (D,E)
$ws.Cells.Item(17, 4).Value = '30.679.85'
$ws.Cells.Item(17, 5).Value = '  +0.39%  '

# Row 18 (D,E)
Set-TextCell 18 4 '13.18'
$ws.Cells.Item(18, 5).Value = '  -1.57%  '

# Row 19 (D,E)
Set-TextCell 19 4 '0.000007521'
$ws.Cells.Item(19, 5).Value = '  +0.27%  '

# Row 20 (E)
$ws.Cells.Item(20, 5).Value = '  +0.14%  '

# Row 21 (D,E)
$ws.Cells.Item(21, 4).Value = '2.143.76'
$ws.Cells.Item(21, 5).Value = '  +1.27%  '

# Row 22 (D,E)
Set-TextCell 22 4 '5.306'
$ws.Cells.Item(22, 5).Value = '  +0.59%  '

# Row 23 (D,E)
Set-TextCell 23 4 '1.002'
$ws.Cells.Item(23, 5).Value = '  +0.06%  '

# Row 24 (D,E)
Set-TextCell 24 4 '6.267'
$ws.Cells.Item(24, 5).Value = '  +1.35%  '

# Row 25 (D,E)
Set-TextCell 25 4 '9.198'
$ws.Cells.Item(25, 5).Value = '  -0.51%  '

# Row 26 (D,E)
Set-TextCell 26 4 '164.50'
$ws.Cells.Item(26, 5).Value = '  +0.53%  '

# Row 27 (D,E)
Set-TextCell 27 4 '18.98'
$ws.Cells.Item(27, 5).Value = '  +1.02%  '

# Row 28 (D,E)
Set-TextCell 28 4 '1.919'
$ws.Cells.Item(28, 5).Value = '  +0.27%  '

# Row 29 (B,C,D,E)
$ws.Cells.Item(29, 2).Value = 'Stellar'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 29 4 '0.09774'
$ws.Cells.Item(29, 5).Value = '  -2.21%  '

# Row 30 (B,C,D,E)
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 30 4 '1.340'
$ws.Cells.Item(30, 5).Value = '  -0.42%  '

# Row 31 (D,E)
Set-TextCell 31 4 '1.485'
$ws.Cells.Item(31, 5).Value = '  -1.44%  '

# Row 32 (D)
Set-TextCell 32 4 '4.313'

# Row 33 (D,E)
Set-TextCell 33 4 '4.163'
$ws.Cells.Item(33, 5).Value = '  +1.25%  '

# Row 34 (E)
$ws.Cells.Item(34, 5).Value = '  +2.54%  '

# Row 35 (D,E)
Set-TextCell 35 4 '1.134'
$ws.Cells.Item(35, 5).Value = '  +1.50%  '

# Row 36 (D,E)
Set-TextCell 36 4 '0.6970'
$ws.Cells.Item(36, 5).Value = '  +0.00%  '

# Row 37 (D,E)
Set-TextCell 37 4 '2.712'
$ws.Cells.Item(37, 5).Value = '  -0.21%  '

# Row 38 (E)
$ws.Cells.Item(38, 5).Value = '  +3.06%  '

# Row 39 (D,E)
Set-TextCell 39 4 '2.845'
$ws.Cells.Item(39, 5).Value = '  +3.40%  '

# Row 40 (D,E)
Set-TextCell 40 4 '75.97'
$ws.Cells.Item(40, 5).Value = '  +3.34%  '

# Row 41 (D,E)
Set-TextCell 41 4 '6.300'
$ws.Cells.Item(41, 5).Value = '  +1.43%  '

# Row 42 (D,E)
Set-TextCell 42 4 '2.008'
$ws.Cells.Item(42, 5).Value = '  +1.91%  '

# Row 43 (D,E)
Set-TextCell 43 4 '0.4280'
$ws.Cells.Item(43, 5).Value = '  +2.54%  '

# Row 44 (E)
$ws.Cells.Item(44, 5).Value = '  -0.01%  '

# Row 45 (D,E)
Set-TextCell 45 4 '0.8305'
$ws.Cells.Item(45, 5).Value = '  -0.31%  '

# Row 46 (D,E)
Set-TextCell 46 4 '101.55'
$ws.Cells.Item(46, 5).Value = '  -0.93%  '

# Row 47 (D)
Set-TextCell 47 4 '9.572'

# Row 48 (B,C,D,E)
$ws.Cells.Item(48, 2).Value = 'Elrond'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 48 4 '35.42'
$ws.Cells.Item(48, 5).Value = '  +0.31%  '

# Row 49 (B,C,D,E)
$ws.Cells.Item(49, 2).Value = 'Aptos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 49 4 '6.996'
$ws.Cells.Item(49, 5).Value = '  +0.62%  '

# Row 50 (D,E)
Set-TextCell 50 4 '912.34'
$ws.Cells.Item(50, 5).Value = '  -1.48%  '

# Row 51 (D,E)
Set-TextCell 51 4 '0.05769'
$ws.Cells.Item(51, 5).Value = '  +2.15%  '
